$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.001.18"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.059.90"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D5").Value = "249.21"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "0.670"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "55.31"
$ws.Range("E8").Value = "  +12.18%  "
$ws.Range("D9").Value = "60.68"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "0.382"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +8.38%  "
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("D13").Value = "15.04"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "2.362.14"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "0.815"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "2.058.31"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "36.934.60"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "0.0₃0944"
$ws.Range("E19").Value = "  +13.57%  "
$ws.Range("D20").Value = "73.48"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  +6.64%  "
$ws.Range("D22").Value = "5.39"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "237.74"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -5.57%  "
$ws.Range("D26").Value = "171.66"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "9.13"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "20.14"
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "4.59"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("D35").Value = "0.0881"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  -6.51%  "
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").Value = "  +24.10%  "
$ws.Range("D41").Value = "18.01"
$ws.Range("E41").Value = "  +7.68%  "
$ws.Range("D42").Value = "0.0225"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "96.98"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "4.13"
$ws.Range("E46").Value = "  +40.98%  "
$ws.Range("E47").Value = "  -49.58%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "4.35"
$ws.Range("E48").Value = "  +11.24%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  +6.43%  "
$ws.Range("D50").Value = "1.301.26"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "2.91"
$ws.Range("E51").Value = "  +0.70%  "
